$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Requested Delivery Date"
$ws.Range("A4").Select()
